$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new "2022-Q4" worksheet, positioned right after "总计"
#    (i.e. right before "2022-Q3"). The easiest way to get an identical
#    layout/formatting to its siblings is to clone the "2022-Q3" sheet and
#    then edit its contents in place.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$srcQ3      = $wb.Worksheets.Item("2022-Q3")

$srcQ3.Copy($null, $totalSheet)      # new copy is placed right after 总计
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The template ("2022-Q3") holds 3 funds; "2022-Q4" only reports 1, so drop
# the two extra data rows (rows 3 and 4), leaving header + single data row.
$newSheet.Rows("3:4").Delete()

# Fund code/name (A2/B2/C2) for the remaining row are already correct
# (same single holding, 159617 / 华夏中证智选500价值稳健策略ETF) -- only the
# figures for the quarter need updating.
$newSheet.Range("D2:G2").NumberFormat = "@"
$newSheet.Cells.Item(2, 4).Value = "1.74"
$newSheet.Cells.Item(2, 5).Value = "97.09"
$newSheet.Cells.Item(2, 6).Value = "1.40"
$newSheet.Cells.Item(2, 7).Value = "0.0244"
$newSheet.Cells.Item(2, 8).Value = 8

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: shift the existing 4 data rows down by
#    one and insert the new 2022-Q4 totals at the top of the data (row 2).
# ---------------------------------------------------------------------------
$ws = $totalSheet

for ($r = 5; $r -ge 2; $r--) {
  $src = $ws.Range("A" + $r + ":D" + $r)
  $dst = $ws.Range("A" + ($r + 1) + ":D" + ($r + 1))
  $src.Copy($dst)
}

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "2022-Q4"
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 0.02

# Renumber the running index in column A for the rows that got shifted down.
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(6, 1).Value = 4

# Restore "总计" as the active sheet/selection, matching the original workbook
# (sheet insert/copy operations above shift the active tab as a side effect).
$totalSheet.Activate() | Out-Null
$totalSheet.Range("A1").Select() | Out-Null
